# Apply the changes described by the diff:
#  1. Update the cached "datetimeFigureOut" footer date from 2021/11/6
#     to 2022/8/5 on the slide master and on every slide layout.
#  2. Update the verse-numbering textboxes on slides 2,3,5,6,8,9,11,12
#     from "( N )" to "( N / 4 )".

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Type -eq 14) {
            $isDate = $false
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    $isDate = $true
                }
            } catch {
            }
            if ($isDate -and $sh.HasTextFrame) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# 1a. Slide master footer date.
Set-DatePlaceholderText $p.SlideMaster.Shapes "2022/8/5"

# 1b. Every slide layout's footer date.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Set-DatePlaceholderText $layout.Shapes "2022/8/5"
}

# 2. Verse-numbering textboxes: "( N )" -> "( N / 4 )".
$numberedSlides = @{
    2  = "( 1 / 4 )"
    3  = "( 1 / 4 )"
    5  = "( 2 / 4 )"
    6  = "( 2 / 4 )"
    8  = "( 3 / 4 )"
    9  = "( 3 / 4 )"
    11 = "( 4 / 4 )"
    12 = "( 4 / 4 )"
}

foreach ($slideIdx in $numberedSlides.Keys) {
    $s = $p.Slides.Item($slideIdx)
    $sh = $s.Shapes.Item($s.Shapes.Count)
    $sh.TextFrame.TextRange.Text = $numberedSlides[$slideIdx]
}
